# Update "想去人数" (column F) values on the 展览 / 本地生活 / 全部类型 sheets
# to match the regenerated data output, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 7747
$ws1.Range("F6").Value  = 99
$ws1.Range("F8").Value  = 2108
$ws1.Range("F9").Value  = 8572
$ws1.Range("F15").Value = 2671
$ws1.Range("F16").Value = 1172
$ws1.Range("F23").Value = 3716
$ws1.Range("F25").Value = 47
$ws1.Range("F26").Value = 37
$ws1.Range("F28").Value = 3216
$ws1.Range("F30").Value = 291
$ws1.Range("F33").Value = 143
$ws1.Range("F34").Value = 347
$ws1.Range("F35").Value = 1020
$ws1.Range("F36").Value = 681
$ws1.Range("F39").Value = 2690
$ws1.Range("F43").Value = 3233
$ws1.Range("F47").Value = 34

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1346

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1346
$ws4.Range("F5").Value  = 7747
$ws4.Range("F6").Value  = 99
$ws4.Range("F8").Value  = 2108
$ws4.Range("F9").Value  = 8572
$ws4.Range("F14").Value = 2671
$ws4.Range("F15").Value = 1172
$ws4.Range("F24").Value = 3716
$ws4.Range("F26").Value = 47
$ws4.Range("F27").Value = 37
$ws4.Range("F29").Value = 3216
$ws4.Range("F32").Value = 143
$ws4.Range("F33").Value = 347
$ws4.Range("F35").Value = 1020
$ws4.Range("F36").Value = 681
$ws4.Range("F40").Value = 2690
$ws4.Range("F44").Value = 3233
$ws4.Range("F47").Value = 34
